{"js": "// Update the answer cells in the \"two-digit \u00f7 one-digit\" worksheet table.\n// Each entry is [oldText, newText] for an exact-match cell (w:t run) replace.\nconst replacements = [\n  [\"58\u00f74=14, 2\", \"89\u00f73=29, 2\"],\n  [\"94\u00f73=31, 1\", \"96\u00f72=48, 0\"],\n  [\"87\u00f78=10, 7\", \"32\u00f78=4, 0\"],\n  [\"35\u00f73=11, 2\", \"18\u00f77=2, 4\"],\n  [\"96\u00f79=10, 6\", \"49\u00f75=9, 4\"],\n  [\"20\u00f76=3, 2\", \"43\u00f73=14, 1\"],\n  [\"99\u00f76=16, 3\", \"84\u00f76=14, 0\"],\n  [\"70\u00f78=8, 6\", \"69\u00f74=17, 1\"],\n  [\"74\u00f75=14, 4\", \"57\u00f72=28, 1\"],\n  [\"94\u00f76=15, 4\", \"98\u00f73=32, 2\"],\n  [\"23\u00f77=3, 2\", \"88\u00f73=29, 1\"],\n  [\"34\u00f77=4, 6\", \"27\u00f76=4, 3\"],\n  [\"91\u00f73=30, 1\", \"87\u00f74=21, 3\"],\n  [\"71\u00f78=8, 7\", \"98\u00f73=32, 2\"],\n  [\"22\u00f76=3, 4\", \"94\u00f72=47, 0\"],\n  [\"84\u00f79=9, 3\", \"89\u00f73=29, 2\"],\n  [\"38\u00f75=7, 3\", \"76\u00f74=19, 0\"],\n  [\"25\u00f76=4, 1\", \"48\u00f73=16, 0\"],\n  [\"23\u00f75=4, 3\", \"44\u00f76=7, 2\"],\n  [\"74\u00f72=37, 0\", \"10\u00f79=1, 1\"],\n  [\"43\u00f74=10, 3\", \"39\u00f76=6, 3\"],\n  [\"67\u00f74=16, 3\", \"98\u00f79=10, 8\"],\n  [\"80\u00f72=40, 0\", \"46\u00f76=7, 4\"],\n  [\"36\u00f79=4, 0\", \"51\u00f77=7, 2\"],\n  [\"53\u00f73=17, 2\", \"45\u00f77=6, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the answer cells in the \"two-digit \u00f7 one-digit\" worksheet table.\n# Each entry is the exact old cell text and its replacement new text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"58\u00f74=14, 2\", \"89\u00f73=29, 2\"),\n    @(\"94\u00f73=31, 1\", \"96\u00f72=48, 0\"),\n    @(\"87\u00f78=10, 7\", \"32\u00f78=4, 0\"),\n    @(\"35\u00f73=11, 2\", \"18\u00f77=2, 4\"),\n    @(\"96\u00f79=10, 6\", \"49\u00f75=9, 4\"),\n    @(\"20\u00f76=3, 2\", \"43\u00f73=14, 1\"),\n    @(\"99\u00f76=16, 3\", \"84\u00f76=14, 0\"),\n    @(\"70\u00f78=8, 6\", \"69\u00f74=17, 1\"),\n    @(\"74\u00f75=14, 4\", \"57\u00f72=28, 1\"),\n    @(\"94\u00f76=15, 4\", \"98\u00f73=32, 2\"),\n    @(\"23\u00f77=3, 2\", \"88\u00f73=29, 1\"),\n    @(\"34\u00f77=4, 6\", \"27\u00f76=4, 3\"),\n    @(\"91\u00f73=30, 1\", \"87\u00f74=21, 3\"),\n    @(\"71\u00f78=8, 7\", \"98\u00f73=32, 2\"),\n    @(\"22\u00f76=3, 4\", \"94\u00f72=47, 0\"),\n    @(\"84\u00f79=9, 3\", \"89\u00f73=29, 2\"),\n    @(\"38\u00f75=7, 3\", \"76\u00f74=19, 0\"),\n    @(\"25\u00f76=4, 1\", \"48\u00f73=16, 0\"),\n    @(\"23\u00f75=4, 3\", \"44\u00f76=7, 2\"),\n    @(\"74\u00f72=37, 0\", \"10\u00f79=1, 1\"),\n    @(\"43\u00f74=10, 3\", \"39\u00f76=6, 3\"),\n    @(\"67\u00f74=16, 3\", \"98\u00f79=10, 8\"),\n    @(\"80\u00f72=40, 0\", \"46\u00f76=7, 4\"),\n    @(\"36\u00f79=4, 0\", \"51\u00f77=7, 2\"),\n    @(\"53\u00f73=17, 2\", \"45\u00f77=6, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 2)\n}\n"}
